$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.523.68'
$ws.Range("E2").Value = '  +0.73%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.66'
$ws.Range("E3").Value = '  +0.07%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.51'
$ws.Range("E5").Value = '  +0.27%  '

$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5076'
$ws.Range("E7").Value = '  -0.85%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3896'
$ws.Range("E8").Value = '  -0.68%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08354'
$ws.Range("E9").Value = '  +0.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.103'
$ws.Range("E10").Value = '  -1.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.78'
$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.219'
$ws.Range("E12").Value = '  -0.50%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.870.87'
$ws.Range("E13").Value = '  -0.68%  '

$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.235'
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.009'
$ws.Range("E16").Value = '  -0.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001103'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.15'
$ws.Range("E18").Value = '  -0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06701'
$ws.Range("E19").Value = '  +0.03%  '

$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("E21").Value = '  -0.23%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.930'
$ws.Range("E22").Value = '  -0.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.558.15'
$ws.Range("E23").Value = '  +0.70%  '

$ws.Range("E24").Value = '  -0.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.232'
$ws.Range("E25").Value = '  -1.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.084.50'
$ws.Range("E26").Value = '  -0.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.70'
$ws.Range("E27").Value = '  +0.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.63'
$ws.Range("E28").Value = '  +0.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.355'
$ws.Range("E29").Value = '  -3.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.10'
$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("E31").Value = '  -1.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.042'
$ws.Range("E32").Value = '  +0.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.792'
$ws.Range("E33").Value = '  -1.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.607'
$ws.Range("E34").Value = '  -0.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02451'
$ws.Range("E35").Value = '  +0.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06548'
$ws.Range("E36").Value = '  +0.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2161'
$ws.Range("E37").Value = '  -0.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.876'
$ws.Range("E38").Value = '  -3.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.059'
$ws.Range("E39").Value = '  +1.93%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.254'
$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.192'
$ws.Range("E41").Value = '  +0.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6425'
$ws.Range("E42").Value = '  -0.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.11'
$ws.Range("E43").Value = '  -0.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.007'
$ws.Range("E44").Value = '  -0.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6035'
$ws.Range("E45").Value = '  -0.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.97'
$ws.Range("E46").Value = '  -0.79%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.687'
$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.011'
$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.216'
$ws.Range("E49").Value = '  +0.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.04'
$ws.Range("E50").Value = '  +0.24%  '

$ws.Range("E51").Value = '  -8.70%  '
